$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the mutable columns (Fecha, Calidad, Volumen, Precio minimo/maximo/promedio, Precio $/Kg)
# for rows 2-41 before any writes, since the edit permutes these values across rows.
$snapshot = @{}
for ($r = 2; $r -le 41; $r++) {
    $row = @{}
    $row["D"] = $ws.Cells.Item($r, 4).Value2
    $row["L"] = $ws.Cells.Item($r, 12).Value2
    $row["M"] = $ws.Cells.Item($r, 13).Value2
    $row["N"] = $ws.Cells.Item($r, 14).Value2
    $row["O"] = $ws.Cells.Item($r, 15).Value2
    $row["P"] = $ws.Cells.Item($r, 16).Value2
    $row["S"] = $ws.Cells.Item($r, 19).Value2
    $snapshot[$r] = $row
}

# Map: destination row -> source row (values originally at the source row now belong at the destination row)
$rowMap = @{
    2 = 25;
    3 = 24;
    4 = 23;
    5 = 19;
    6 = 28;
    7 = 27;
    8 = 22;
    9 = 26;
    10 = 36;
    11 = 34;
    12 = 33;
    13 = 7;
    14 = 12;
    15 = 13;
    16 = 17;
    17 = 5;
    18 = 38;
    19 = 39;
    20 = 10;
    21 = 41;
    22 = 6;
    23 = 14;
    24 = 15;
    25 = 29;
    26 = 31;
    27 = 30;
    28 = 18;
    29 = 37;
    30 = 4;
    31 = 20;
    32 = 40;
    33 = 2;
    34 = 3;
    35 = 21;
    36 = 32;
    37 = 35;
    38 = 8;
    39 = 9;
    40 = 16;
    41 = 11
}

foreach ($dst in $rowMap.Keys) {
    $src = $rowMap[$dst]
    $data = $snapshot[$src]
    $ws.Cells.Item($dst, 4).Value2 = $data["D"]
    $ws.Cells.Item($dst, 12).Value2 = $data["L"]
    $ws.Cells.Item($dst, 13).Value2 = $data["M"]
    $ws.Cells.Item($dst, 14).Value2 = $data["N"]
    $ws.Cells.Item($dst, 15).Value2 = $data["O"]
    $ws.Cells.Item($dst, 16).Value2 = $data["P"]
    $ws.Cells.Item($dst, 19).Value2 = $data["S"]
}